$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 17:10:20"
$wsZhCn.Range("H2").Value = "2016-03-19 17:11:03"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 17:10:28"
$wsDeDe.Range("H2").Value = "2016-03-19 17:11:16"
